# Adds an "Image" column (D) with product image URLs, sourced from the
# menu photo hosting links. D2 is turned into a live hyperlink (matching
# the Hyperlink cell style Excel auto-creates); D3:D35 are plain text URLs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$imageUrls = @(
    @(2, "https://i.ibb.co/DKRtvJ7/800-6012d8a3a0807.jpg"),
    @(3, "https://i.ibb.co/4FddBW5/800-6012d8e670df3.jpg"),
    @(4, "https://i.ibb.co/2t7zHgL/800-6012d82e751c0.jpg"),
    @(5, "https://i.ibb.co/0QBLMr3/800-6012d9626749b.jpg"),
    @(6, "https://i.ibb.co/1XXDSmb/800-66392c99b094d.png"),
    @(7, "https://i.ibb.co/1rDv91n/800-66392dae99407.png"),
    @(8, "https://i.ibb.co/hBsYd8m/800-6012d1d6f1af8.jpg"),
    @(9, "https://i.ibb.co/P90SVMs/800-6012d21eefaf5.jpg"),
    @(10, "https://i.ibb.co/hBsYd8m/800-6012d1d6f1af8.jpg"),
    @(11, "https://i.ibb.co/hBsYd8m/800-6012d1d6f1af8.jpg"),
    @(12, "https://i.ibb.co/hBsYd8m/800-6012d1d6f1af8.jpg"),
    @(13, "https://i.ibb.co/PZLbXzg/800-6012d13256a8b.jpg"),
    @(14, "https://i.ibb.co/3TJGtbR/800-60134f31a1f32.jpg"),
    @(15, "https://i.ibb.co/3TJGtbR/800-60134f31a1f32.jpg"),
    @(16, "https://i.ibb.co/JmRZV7s/800-6012aab2c15a7.jpg"),
    @(17, "https://i.ibb.co/PZG8Njq/800-6012aa1e6fa36.jpg"),
    @(18, "https://i.ibb.co/PDgF2rJ/800-6012aa618cff2.jpg"),
    @(19, "https://i.ibb.co/JmRZV7s/800-6012aab2c15a7.jpg"),
    @(20, "https://i.ibb.co/zFxv5Ds/800-6012cf167d3bb.jpg"),
    @(21, "https://i.ibb.co/r01ytKq/800-6012ce8dd3461.jpg"),
    @(22, "https://i.ibb.co/t481pHP/800-6012ce5b80fc7.jpg"),
    @(23, "https://i.ibb.co/LPzxzRR/800-6012ce123334f.jpg"),
    @(24, "https://i.ibb.co/NL8T1TR/800-6012adc541a5b.jpg"),
    @(25, "https://i.ibb.co/Gd5PdtL/800-6012ae325a211.jpg"),
    @(26, "https://i.ibb.co/kDBKRbS/800-6012acac63175.jpg"),
    @(27, "https://i.ibb.co/h2t8vKN/800-6012ad30ba62d.jpg"),
    @(28, "https://i.ibb.co/bRrt19b/800-6012ad821225e.jpg"),
    @(29, "https://i.ibb.co/nRpWWtc/800-6012ac57acc78.jpg"),
    @(30, "https://i.ibb.co/fFJwgXM/800-6012a8d2ca653.jpg"),
    @(31, "https://i.ibb.co/KNkJ5nQ/800-6012a96d269da.jpg"),
    @(32, "https://i.ibb.co/kDBKRbS/800-6012acac63175.jpg"),
    @(33, "https://i.ibb.co/KjyjT6F/800-6012a8655da02.jpg"),
    @(34, "https://i.ibb.co/nktXNQ3/800-6012a812af173.jpg"),
    @(35, "https://i.ibb.co/LdKwhWm/800-6012a6070a898.jpg")
)

foreach ($pair in $imageUrls) {
    $row = $pair[0]
    $url = $pair[1]
    $cell = $ws.Cells.Item($row, 4)
    if ($row -eq 2) {
        $ws.Hyperlinks.Add($cell, $url, "", "", $url)
    } else {
        $cell.Value = $url
    }
}

# Header for the new column, added last so it lands at the end of the
# shared-string table (after the 28 unique image URLs).
$ws.Cells.Item(1, 4).Value = "Image"

# Restore the cursor/selection position recorded in the saved workbook.
$ws.Range("I6").Select()
